$d = $word.ActiveDocument

$replacements = @(
    @{old = "326÷8=40, 6"; new = "447÷4=111, 3"},
    @{old = "481÷2=240, 1"; new = "990÷6=165, 0"},
    @{old = "741÷8=92, 5"; new = "236÷8=29, 4"},
    @{old = "272÷7=38, 6"; new = "493÷6=82, 1"},
    @{old = "529÷8=66, 1"; new = "538÷5=107, 3"},
    @{old = "972÷3=324, 0"; new = "249÷9=27, 6"},
    @{old = "342÷9=38, 0"; new = "373÷3=124, 1"},
    @{old = "978÷3=326, 0"; new = "491÷7=70, 1"},
    @{old = "688÷7=98, 2"; new = "312÷3=104, 0"},
    @{old = "743÷4=185, 3"; new = "467÷2=233, 1"},
    @{old = "341÷3=113, 2"; new = "172÷8=21, 4"},
    @{old = "206÷2=103, 0"; new = "999÷2=499, 1"},
    @{old = "453÷5=90, 3"; new = "937÷7=133, 6"},
    @{old = "631÷6=105, 1"; new = "918÷9=102, 0"},
    @{old = "875÷4=218, 3"; new = "912÷3=304, 0"},
    @{old = "481÷6=80, 1"; new = "288÷5=57, 3"},
    @{old = "602÷4=150, 2"; new = "605÷7=86, 3"},
    @{old = "171÷4=42, 3"; new = "688÷3=229, 1"},
    @{old = "831÷3=277, 0"; new = "298÷5=59, 3"},
    @{old = "218÷8=27, 2"; new = "784÷5=156, 4"},
    @{old = "916÷6=152, 4"; new = "864÷4=216, 0"},
    @{old = "200÷9=22, 2"; new = "795÷6=132, 3"},
    @{old = "895÷8=111, 7"; new = "774÷5=154, 4"},
    @{old = "451÷6=75, 1"; new = "891÷2=445, 1"},
    @{old = "750÷6=125, 0"; new = "938÷5=187, 3"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
